$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 2 keeps its BPID (HT0073); only the Business Model Name text changes.
$ws1.Range("B2").Value = "Expedia Marketing Fee"

# Row 3: BPID HT0090 -> HT0089; name becomes "Expedia Marketing Fee"
$ws1.Range("A3").Value = "HT0089"
$ws1.Range("B3").Value = "Expedia Marketing Fee"

# Row 4: BPID HT0089 -> HT0254; name becomes "Expedia Marketing Fee (Penalty)"
# (set B4 before A4 so the shared-string table allocates new entries in the
# same order as the target workbook: "...(Penalty)" text, then "HT0254")
$ws1.Range("B4").Value = "Expedia Marketing Fee (Penalty)"
$ws1.Range("A4").Value = "HT0254"

# Column A widened (manually resized) and no longer auto "best fit".
$ws1.Columns.Item(1).ColumnWidth = 10.67

# Active selection moved to C2 after the edits.
$ws1.Range("C2").Select()
